$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Row 56: username "retchygod" -> "24feetofsmoke"
$ws.Range("B56").Value = "24feetofsmoke"

# 2. Row 128: name "Наталья" -> "Натали"
$ws.Range("C128").Value = "Натали"

# 3. Insert a new row after row 142 (Сергей Николаевич) for Сергей Цыбулин
$ws.Rows.Item(143).Insert()

# The id column stores purely-numeric ids as text. Stage the id as text in a
# scratch cell (formatted as Text) then paste-by-value into place so the
# destination cell picks up the shared-string/text type without a number
# format being attached to it.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "8628968557"
$ws.Range("Z1").Copy()
$ws.Range("A143").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("B143").Value = "cybulin2018"
$ws.Range("C143").Value = "Сергей Цыбулин"

# 4. Remove the row with chill.resort (now shifted to row 154)
$ws.Rows.Item(154).Delete()
